# Generate Report for Handoff
# Updates the localization-status report: status moves from
# "In Translation" to "Ready for handoff", refreshes the associated
# timestamps, and widens the Status columns to fit the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---------------------------------------------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-02 07:07:46"

$wsOverview.Columns.Item(5).ColumnWidth = 16.3333333333333
$wsOverview.Columns.Item(6).ColumnWidth = 16.3333333333333

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-02 07:07:42"

$wsZhCn.Columns.Item(3).ColumnWidth = 16.3333333333333

# --- de-de sheet --------------------------------------------------------
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-02 07:07:46"

$wsDeDe.Columns.Item(3).ColumnWidth = 16.3333333333333
